$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A8").Value = -21.25890000000002
$ws.Range("A10").Value = -20.51169999999998
$ws.Range("A12").Value = -22.46240000000002
$ws.Range("A18").Value = -22.32760000000002
